# Updating SS concentrations from LISST data
# Fill in the "Time" column (B) on the "Storm 1" sheet with the
# sample-collection times recovered from the LISST instrument log.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Storm 1")
$ws.Activate()

$times = [ordered]@{
    "B2"  = 0.65972222222222221
    "B3"  = 0.67291666666666661
    "B4"  = 0.69444444444444453
    "B5"  = 0.70486111111111116
    "B6"  = 0.80486111111111114
    "B7"  = 0.87638888888888899
    "B10" = 0.62083333333333335
    "B11" = 0.62638888888888888
    "B12" = 0.67569444444444438
    "B13" = 0.70694444444444438
    "B14" = 0.74861111111111101
    "B15" = 0.84236111111111101
}

foreach ($addr in $times.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = $times[$addr]
    $cell.NumberFormat = "h:mm"
}

# Restore the active cell/selection recorded on this sheet.
[void]$ws.Range("K15").Select()
